{"js": "// Remove right-alignment from the date/period cell (the second cell)\n// of every \"Professional Experience\" entry table so its text falls\n// back to the default (left) alignment.\n//\n// Document layout: the first table in the body is the header table\n// (name / contact info) and must stay untouched. Every subsequent\n// top-level table is a single-row, two-column entry: left cell holds\n// the employer name, right cell holds the employment period and is\n// currently forced right-aligned (w:jc val=\"right\"); we clear that.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 1; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  table.rows.load(\"items\");\n  await context.sync();\n\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n    await context.sync();\n\n    const lastCell = row.cells.items[row.cells.items.length - 1];\n    const paragraphs = lastCell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    for (const paragraph of paragraphs.items) {\n      paragraph.load(\"alignment\");\n    }\n    await context.sync();\n\n    for (const paragraph of paragraphs.items) {\n      if (paragraph.alignment === \"Right\") {\n        paragraph.alignment = \"Left\";\n      }\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove right-alignment from the date/period cell (the last cell) of\n# every \"Professional Experience\" entry table so its text falls back\n# to the default (left) alignment.\n#\n# Document layout: Tables.Item(1) is the header table (name / contact\n# info) and must stay untouched. Every subsequent top-level table is a\n# single-row, two-column entry: first cell holds the employer name,\n# last cell holds the employment period and is currently forced\n# right-aligned (w:jc val=\"right\"); we clear that back to the default.\n$d = $word.ActiveDocument\n\nfor ($i = 2; $i -le $d.Tables.Count; $i++) {\n    $table = $d.Tables.Item($i)\n    foreach ($row in $table.Rows) {\n        $lastCell = $row.Cells.Item($row.Cells.Count)\n        foreach ($para in $lastCell.Range.Paragraphs) {\n            if ($para.Alignment -eq 2) {\n                $para.Alignment = 0\n            }\n        }\n    }\n}\n"}
